$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63; this shifts the existing rows 63-149 down to 64-150,
# which matches the bulk of the diff (all of those rows simply retain their previous
# values one row further down).
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new data record.
$ws.Range("A63").Value = 10
$ws.Range("B63").Value = "Vega Modelo de Temuco"
$ws.Range("C63").Value = "La Araucanía"
$ws.Range("D63").Value = 45195
$ws.Range("E63").Value = 9
$ws.Range("F63").Value = 100112010
$ws.Range("G63").Value = "Achicoria"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 80
$ws.Range("K63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = 10000
$ws.Range("N63").Value = "$/caja 18 unidades"
$ws.Range("O63").Value = "Región Metropolitana"
$ws.Range("P63").Value = 556
$ws.Range("Q63").Value = 18
$ws.Range("R63").Value = "Hortaliza"
